$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated s_vals data after regenerating to filter save games.
# Columns: B=TB, C=d2S, D=K, E=IP, F=Win, G=sum

$data = @(
    @{ Row = 2; B = 3.182878228561681;  C = 1.65323645889881;    D = 0.1529057820181812;  E = 0.4998867070740569; G = 5.488907176552729  },
    @{ Row = 3; B = 0.1554434735375247; C = 0.05231270169004087; D = 0.7127328510149897;  E = 0.4998867070740569; G = 1.420375733316612   },
    @{ Row = 4; B = 0.3464964993005633; C = 1.65323645889881;    D = 0.7127328510149897;  E = 6.48142807727062;   G = 9.193893886484982  },
    @{ Row = 5; B = 3.182878228561681;  C = 1.65323645889881;    D = 16.98373111632243;   E = 0.4998867070740569; G = 22.31973251085698  },
    @{ Row = 6; B = 3.182878228561681;  C = 1.65323645889881;    D = 16.98373111632243;   E = 0.4998867070740569; G = 22.31973251085698  },
    @{ Row = 7; B = 3.182878228561681;  C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538  }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 7).Value = $entry.G
}
